$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at 624-625, shifting existing rows 624:708 down to 626:710
$ws.Rows("624:625").Insert()

# Populate new row 624 (Primera) with the new weekly data point
$ws.Cells.Item(624, 1).Value = 9
$ws.Cells.Item(624, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(624, 3).Value = "Metropolitana"
$ws.Cells.Item(624, 4).Value = 44918
$ws.Cells.Item(624, 5).Value = 13
$ws.Cells.Item(624, 6).Value = 100114014
$ws.Cells.Item(624, 7).Value = "Betarraga"
$ws.Cells.Item(624, 8).Value = "Sin especificar"
$ws.Cells.Item(624, 9).Value = "Primera"
$ws.Cells.Item(624, 10).Value = 10600
$ws.Cells.Item(624, 11).Value = 90
$ws.Cells.Item(624, 12).Value = 100
$ws.Cells.Item(624, 13).Value = 95
$ws.Cells.Item(624, 14).Value = "`$/unidad"
$ws.Cells.Item(624, 15).Value = "Región Metropolitana"
$ws.Cells.Item(624, 16).Value = 95
$ws.Cells.Item(624, 17).Value = 1
$ws.Cells.Item(624, 18).Value = "Hortaliza"

# Populate new row 625 (Segunda) with the new weekly data point
$ws.Cells.Item(625, 1).Value = 9
$ws.Cells.Item(625, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(625, 3).Value = "Metropolitana"
$ws.Cells.Item(625, 4).Value = 44918
$ws.Cells.Item(625, 5).Value = 13
$ws.Cells.Item(625, 6).Value = 100114014
$ws.Cells.Item(625, 7).Value = "Betarraga"
$ws.Cells.Item(625, 8).Value = "Sin especificar"
$ws.Cells.Item(625, 9).Value = "Segunda"
$ws.Cells.Item(625, 10).Value = 5200
$ws.Cells.Item(625, 11).Value = 70
$ws.Cells.Item(625, 12).Value = 70
$ws.Cells.Item(625, 13).Value = 70
$ws.Cells.Item(625, 14).Value = "`$/unidad"
$ws.Cells.Item(625, 15).Value = "Región Metropolitana"
$ws.Cells.Item(625, 16).Value = 70
$ws.Cells.Item(625, 17).Value = 1
$ws.Cells.Item(625, 18).Value = "Hortaliza"

Write-Host "Done"
